# New crime data collected - weekly 43rd Precinct CompStat update
# Updates the "Volume/Number" header, the reporting week date range,
# and refreshes the crime-complaint statistics table (rows 15-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings)
# ---------------------------------------------------------------------------

# "Volume 32   Number  5" -> "Volume 32   Number  6"
$a8 = $ws.Range("A8")
$a8txt = $a8.Text
$a8.Characters($a8txt.Length, 1).Text = "6"

# "Report Covering the Week  1/27/2025  Through  2/2/2025"
#   -> "Report Covering the Week  2/3/2025  Through  2/9/2025"
$c9 = $ws.Range("C9")

$c9txt = $c9.Text
$oldStart = "1/27/2025"
$newStart = "2/3/2025"
$idxStart = $c9txt.IndexOf($oldStart)
$c9.Characters($idxStart + 1, $oldStart.Length).Text = $newStart

$c9txt2 = $ws.Range("C9").Text
$oldEnd = "2/2/2025"
$newEnd = "2/9/2025"
$idxEnd = $c9txt2.IndexOf($oldEnd)
$ws.Range("C9").Characters($idxEnd + 1, $oldEnd.Length).Text = $newEnd

# ---------------------------------------------------------------------------
# Helper: cells that flip between a numeric style and the "no data" text
# style (shared strings "0" / "***.*", cell style 13) need their style
# copied from a stable donor cell before/while the value is written, since
# plain .Value assignment keeps the previous style.
# ---------------------------------------------------------------------------

# Donor for text "0" (shared string 20) with style 13
$donorText0 = $ws.Range("D15")
# Donor for text "***.*" (shared string 21) with style 13
$donorTextDashes = $ws.Range("E15")
# Donor for numeric style 14 with value 1 (stays constant across the edit)
$donorNum1 = $ws.Range("F14")
# Donor for numeric style 15 with value 0 (stays constant across the edit)
$donorPct0 = $ws.Range("H15")

# Row 15 (Rape): C goes from a number to the "0" placeholder text
$donorText0.Copy($ws.Range("C15"))
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -50

# Row 16 (Robbery)
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 32
$ws.Range("H16").Value = 3.225806451612
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -2.040816326530
$ws.Range("L16").Value = -35.135135135135
$ws.Range("M16").Value = 33.333333333333
$ws.Range("N16").Value = -78.082191780821

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 38.461538461538
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 30.952380952381
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 63
$ws.Range("K17").Value = 4.761904761904
$ws.Range("L17").Value = -35.294117647058
$ws.Range("M17").Value = 40.425531914893
$ws.Range("N17").Value = -26.666666666666

# Row 18 (Burglary)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -35.483870967741
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -26.530612244898
$ws.Range("M18").Value = -21.739130434782
$ws.Range("N18").Value = -77.358490566037

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 78
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -7.142857142857
$ws.Range("I19").Value = 97
$ws.Range("J19").Value = 111
$ws.Range("K19").Value = -12.612612612612
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 90.196078431372
$ws.Range("N19").Value = 46.969696969697

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -15.384615384615
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = -19.047619047619
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = -12.5
$ws.Range("L20").Value = -43.678160919540
$ws.Range("M20").Value = 133.333333333333
$ws.Range("N20").Value = -79.835390946502

# Row 21 (TOTAL)
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -13.432835820895
$ws.Range("F21").Value = 223
$ws.Range("G21").Value = 233
$ws.Range("H21").Value = -4.291845493562
$ws.Range("I21").Value = 302
$ws.Range("J21").Value = 318
$ws.Range("K21").Value = -5.031446540880
$ws.Range("L21").Value = -26.876513317191
$ws.Range("M21").Value = 47.317073170731
$ws.Range("N21").Value = -61.916771752837

# Row 23 (Housing)
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -16.666666666666
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = -10
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = -29.032258064516
$ws.Range("L23").Value = -43.589743589743
$ws.Range("M23").Value = 29.411764705882

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 57.894736842105
$ws.Range("F24").Value = 150
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = 12.781954887218
$ws.Range("I24").Value = 213
$ws.Range("J24").Value = 191
$ws.Range("K24").Value = 11.518324607329
$ws.Range("L24").Value = -6.167400881057
$ws.Range("M24").Value = 45.890410958904

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 4.761904761904
$ws.Range("I25").Value = 61
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = -6.153846153846
$ws.Range("L25").Value = -44.036697247706

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -8.333333333333
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 92
$ws.Range("H26").Value = -30.434782608695
$ws.Range("I26").Value = 99
$ws.Range("J26").Value = 124
$ws.Range("K26").Value = -20.161290322580
$ws.Range("L26").Value = -13.157894736842
$ws.Range("M26").Value = -37.735849056603

# Row 27 (UCR Rape*): C, D, E switch to the "0"/"***.*" placeholder text
$donorText0.Copy($ws.Range("C27"))
$donorText0.Copy($ws.Range("D27"))
$donorTextDashes.Copy($ws.Range("E27"))
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25

# Row 28 (Other Sex Crimes)
$ws.Range("D28").Value = 3
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -78.571428571428
$ws.Range("L28").Value = -66.666666666666

# Row 29 (Shooting Vic.): D, E switch to "0"/"***.*" placeholder text
$donorText0.Copy($ws.Range("D29"))
$donorTextDashes.Copy($ws.Range("E29"))
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("L29").Value = 0

# Row 30 (Shooting Inc.): D, E switch to "0"/"***.*" placeholder text
$donorText0.Copy($ws.Range("D30"))
$donorTextDashes.Copy($ws.Range("E30"))
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("L30").Value = 0

# Row 33 (Traffic Fatalities): C-K switch from placeholder text to real
# numeric data (1,1,0% change repeated three times)
$donorNum1.Copy($ws.Range("C33"))
$donorNum1.Copy($ws.Range("D33"))
$donorPct0.Copy($ws.Range("E33"))
$donorNum1.Copy($ws.Range("F33"))
$donorNum1.Copy($ws.Range("G33"))
$donorPct0.Copy($ws.Range("H33"))
$donorNum1.Copy($ws.Range("I33"))
$donorNum1.Copy($ws.Range("J33"))
$donorPct0.Copy($ws.Range("K33"))
